$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DOE")

# Row 2 (Evaporator Temperature): switch from Continuous (Mean/StdDev) to Discrete (Max/Min/Step)
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 30
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0.1
$ws.Range("H2").Value = "Discrete"

# Row 3 (Condenser Temperature): update Max/Min
$ws.Range("D3").Value = 65
$ws.Range("E3").Value = 35

# Row 5 (Capacity): update Max/Min/Step
$ws.Range("D5").Value = 20000
$ws.Range("E5").Value = 5000
$ws.Range("F5").Value = 100
